# Fix bug in genSymbolValsWithPos: the rows in the "wrss" test sheet were
# generated/written in the wrong order. This reorders the data rows (A:F,
# rows 3-23) back to the correct row positions while preserving all of the
# underlying values (symbol id + the 5 reel counts) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is: row number, then the values for columns A..F that must
# end up on that row.
$rows = @(
    @(3, 201,9,30,15,45,30),
    @(4, 1201,2,10,10,10,10),
    @(5, 1203,3,15,15,15,15),
    @(6, 501,9,52,30,75,45),
    @(7, 701,3,90,45,97,15),
    @(8, 601,9,60,67,60,42),
    @(9, 101,9,30,15,60,15),
    @(10, 902,1,0,0,0,0),
    @(11, 401,9,48,67,75,45),
    @(12, 1202,2,10,10,10,10),
    @(13, 1001,18,30,75,60,72),
    @(14, 301,6,45,30,60,45),
    @(15, 801,3,67,65,52,45),
    @(16, 502,0,4,0,0,0),
    @(17, 802,0,4,5,4,0),
    @(18, 1101,0,15,30,30,0),
    @(19, 1,0,2,2,2,2),
    @(20, 2,0,2,2,2,2),
    @(21, 3,0,3,3,3,3),
    @(22, 602,0,0,4,0,9),
    @(23, 402,0,0,4,0,0)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $entry[$col]
    }
}
